$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADC_100MS")

# Update existing cell D24: 110 -> 210
$ws.Range("D24").Value = 210

# Add new row 25 with instance "mem_i2"
$ws.Range("B25").Value = "mem_i2"
$ws.Range("C25").Value = "MEM1"
$ws.Range("D25").Value = 210
$ws.Range("E25").Value = 160

# Match formatting of the row above (D/E cells centered, matching style of D23/E-column default)
$ws.Range("D25").HorizontalAlignment = -4108
$ws.Range("E25").HorizontalAlignment = -4108

# Update selection to new active cell
$ws.Range("E25").Select()
